# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @(1.505614041169197, 87981.0709163148, 157.8057217802531, 246.9852506941017, 88387.36750283033)
    3 = @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.705647867635037)
    4 = @(1.505614041169197, 1.65323645889881, 16.98373111632243, 6.48142807727062, 26.62400969366105)
    5 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 12.0302756157461)
    6 = @(0.1554434735375247, 9.226618575922256, 157.8057217802531, 6.48142807727062, 173.6692119069835)
    7 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    8 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 10.35301142835362)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
